# Implement the "ameliorations" level-2 (apprentice 1) game-logic row,
# matching the author's commit "gamelogic almost fully implemented".

$wb = $excel.ActiveWorkbook

# --- ameliorations sheet: fill in row 4 (ID 2, "apprentice 1") ---
$wsAme = $wb.Worksheets.Item("ameliorations")
$wsAme.Activate()

$wsAme.Range("D4").Value = 200
$wsAme.Range("E4").Value = 3
$wsAme.Range("F4").Value = 3
$wsAme.Range("H4").Formula = "=D4*E4^F4"

# matches the new selection recorded in the sheet's saved view
$wsAme.Range("F4").Select()

# --- techs sheet: just a cursor/selection move ---
$wsTechs = $wb.Worksheets.Item("techs")
$wsTechs.Activate()
$wsTechs.Range("D5").Select()
